$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '61.926.59'
$c.Style = 'Normal'
$ws.Range('E2').Value = '  -0.77%  '

# Row 3
$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '3.411.83'
$c.Style = 'Normal'
$ws.Range('E3').Value = '  -0.59%  '

# Row 4
$ws.Range('E4').Value = '  -0.10%  '

# Row 5
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '410.08'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +0.66%  '

# Row 6
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '129.15'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  -0.83%  '

# Row 7
$ws.Range('E7').Value = '  +6.14%  '

# Row 8
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = 'Normal'
$ws.Range('E8').Value = '  +0.04%  '

# Row 9
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.737'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  +5.78%  '

# Row 10
$ws.Range('E10').Value = '  -1.04%  '

# Row 11
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '43.57'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  +3.30%  '

# Row 12
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '9.37'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  +9.75%  '

# Row 13
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '0.0000220'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  +30.49%  '

# Row 14
$ws.Range('E14').Value = '  -0.37%  '

# Row 15
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '21.39'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  +7.30%  '

# Row 16
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '3.952.64'
$c.Style = 'Normal'
$ws.Range('E16').Value = '  -0.82%  '

# Row 17
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '3.420.65'
$c.Style = 'Normal'
$ws.Range('E17').Value = '  -0.02%  '

# Row 18
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '12.54'
$c.Style = 'Normal'
$ws.Range('E18').Value = '  +8.39%  '

# Row 19
$ws.Range('E19').Value = '  +6.45%  '

# Row 20
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '61.948.39'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  -0.84%  '

# Row 21
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '451.19'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  +43.35%  '

# Row 22
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '91.68'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  +8.04%  '

# Row 23
$ws.Range('E23').Value = '  +0.32%  '

# Row 24
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '13.24'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  +2.91%  '

# Row 25
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '3.31'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  +3.95%  '

# Row 26
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '9.36'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  +14.35%  '

# Row 27
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '33.22'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  +10.87%  '

# Row 28
$ws.Range('E28').Value = '  +0.71%  '

# Row 29
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '7.69'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  -1.19%  '

# Row 30
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '2.71'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  +0.09%  '

# Row 31
$ws.Range('E31').Value = '  +5.20%  '

# Row 32
$ws.Range('E32').Value = '  -2.04%  '

# Row 33
$ws.Range('E33').Value = '  -0.38%  '

# Row 34
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '42.64'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  -3.59%  '

# Row 35
$ws.Range('E35').Value = '  +0.07%  '

# Row 36
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '0.0507'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  +4.35%  '

# Row 37
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '53.81'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  +4.44%  '

# Row 38
$ws.Range('E38').Value = '  -0.17%  '

# Row 39
$ws.Range('B39').Value = 'Stellar'
$ws.Range('C39').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '0.137'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  +8.54%  '

# Row 40
$ws.Range('B40').Value = 'LidoDAOToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '3.40'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  +2.04%  '

# Row 41
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '2.96'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  -0.51%  '

# Row 42
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '0.318'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  -0.39%  '

# Row 43
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '143.92'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  -0.02%  '

# Row 44
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '4.35'
$c.Style = 'Normal'

# Row 45
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '2.01'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  +1.11%  '

# Row 46
$ws.Range('E46').Value = '  +15.16%  '

# Row 47
$ws.Range('E47').Value = '  -1.51%  '

# Row 48
$ws.Range('E48').Value = '  +22.98%  '

# Row 49
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '22.54'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  +5.57%  '

# Row 50
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '2.17'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  +7.12%  '

# Row 51
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '3.753.83'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  -0.62%  '
